$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values for columns B..Q (the same set of values is written to
# every data row, 2..26 - matching the diff which updates all rows identically).
$values = @(
    0.999988373277124198423848611128,
    0.999076334919017483571224147454,
    0.999999999999836242103867789410,
    0.999999940172837220586643525166,
    0.999999982861204173012481533078,
    0.000010853035450360880667143582,
    0.000862200808886026985253248078,
    0.000000000000147085511348572506,
    0.000000027346381937709721326107,
    0.000000013673264511610529775351,
    0.000180565691205387295666801029,
    0.003294394549892420195197351163,
    0.999906986216993587390788889024,
    0.003434643746218192950775405237,
    64.862131504189036945717816706747,
    90.458523826421242119977250695229
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
